$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test time header
$ws.Range("A1").Value = "测试时间:  : 2025-11-07 13:54:00"

# Update model name and existing metric columns for row 5
$ws.Range("B5").Value = "Qwen3-4B-Instruct-2507"
$ws.Range("C5").Value = "0.3871 ± 0.0049"
$ws.Range("D5").Value = "0.4345 ± 0.0049"
$ws.Range("E5").Value = "241.51S"

# Apply the same style as the existing row-5 cells to the newly added cells
$ws.Range("E5").Copy()
$ws.Range("F5:N5").PasteSpecial(-4122)

# Fill in the newly added columns (vllm bfloat16 / FP8 / AWQ results)
$ws.Range("F5").Value = "0.3878 ± 0.0049"
$ws.Range("G5").Value = "0.4352 ± 0.0049"
$ws.Range("H5").Value = "293.21S"
$ws.Range("I5").Value = "0.3889 ± 0.0049"
$ws.Range("J5").Value = "0.4362 ± 0.0049"
$ws.Range("K5").Value = "236.85S"
$ws.Range("L5").Value = "0.3878 ± 0.0049"
$ws.Range("M5").Value = "0.4349 ± 0.0049"
$ws.Range("N5").Value = "275.39S"

$excel.CutCopyMode = $false
